$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = "'52.219.76"
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Formula = "'2.822.59"
$ws.Range('E3').Value = '  +1.16%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Formula = "'356.64"
$ws.Range('E5').Value = '  +2.75%  '
$ws.Range('D6').Formula = "'112.02"
$ws.Range('E6').Value = '  -3.31%  '
$ws.Range('D7').Formula = "'0.571"
$ws.Range('E7').Value = '  +4.00%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Formula = "'0.599"
$ws.Range('E9').Value = '  +1.99%  '
$ws.Range('D10').Formula = "'41.01"
$ws.Range('E10').Value = '  -4.20%  '
$ws.Range('E11').Value = '  +0.72%  '
$ws.Range('E12').Value = '  +1.02%  '
$ws.Range('D13').Formula = "'19.88"
$ws.Range('E13').Value = '  -0.74%  '
$ws.Range('D14').Formula = "'7.77"
$ws.Range('E14').Value = '  -1.40%  '
$ws.Range('D15').Formula = "'3.266.95"
$ws.Range('E15').Value = '  +1.19%  '
$ws.Range('D16').Formula = "'2.825.17"
$ws.Range('E16').Value = '  +1.55%  '
$ws.Range('D17').Formula = "'0.928"
$ws.Range('E17').Value = '  +4.42%  '
$ws.Range('D18').Formula = "'52.100.26"
$ws.Range('E18').Value = '  -0.03%  '
$ws.Range('D19').Formula = "'7.52"
$ws.Range('E19').Value = '  +4.06%  '
$ws.Range('E20').Value = '  -0.47%  '
$ws.Range('D21').Formula = "'13.52"
$ws.Range('E21').Value = '  +0.81%  '
$ws.Range('D22').Formula = "'0.0₃0997"
$ws.Range('E22').Value = '  +2.03%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').Formula = "'271.55"
$ws.Range('E23').Value = '  +0.63%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Formula = "'70.49"
$ws.Range('E24').Value = '  +0.82%  '
$ws.Range('E25').Value = '  +1.81%  '
$ws.Range('D26').Formula = "'26.94"
$ws.Range('E26').Value = '  +1.12%  '
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('E28').Value = '  +1.24%  '
$ws.Range('E29').Value = '  +0.51%  '
$ws.Range('E30').Value = '  +8.29%  '
$ws.Range('E31').Value = '  +2.55%  '
$ws.Range('B32').Value = 'InjectiveProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D32').Formula = "'35.21"
$ws.Range('E32').Value = '  +1.28%  '
$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D33').Formula = "'52.44"
$ws.Range('E33').Value = '  +4.47%  '
$ws.Range('D34').Formula = "'5.93"
$ws.Range('E34').Value = '  +4.08%  '
$ws.Range('E35').Value = '  +12.71%  '
$ws.Range('E36').Value = '  +3.85%  '
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('E38').Value = '  +2.04%  '
$ws.Range('D40').Formula = "'18.43"
$ws.Range('E40').Value = '  -1.54%  '
$ws.Range('E41').Value = '  +1.79%  '
$ws.Range('D42').Formula = "'127.41"
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('E43').Value = '  -3.51%  '
$ws.Range('D44').Formula = "'23.16"
$ws.Range('E44').Value = '  -0.56%  '
$ws.Range('E45').Value = '  -1.07%  '
$ws.Range('D46').Formula = "'2.092.17"
$ws.Range('E46').Value = '  +1.37%  '
$ws.Range('D47').Formula = "'3.36"
$ws.Range('E47').Value = '  +0.92%  '
$ws.Range('D48').Formula = "'2.27"
$ws.Range('E48').Value = '  -2.40%  '
$ws.Range('E49').Value = '  +7.34%  '
$ws.Range('D50').Formula = "'0.966"
$ws.Range('E50').Value = '  +0.20%  '
$ws.Range('E51').Value = '  +3.04%  '
